$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 11. This shifts the previous rows 11-19
# down to 12-20, preserving their existing values untouched.
$ws.Rows.Item(11).Insert()

# --- Row 10: this record (doc 25018849 / NORA MARIA BLANDON SERNA) is now
# flagged as existing only on the "Instructores" side - clear the Sofia-side
# and the leading descriptive columns, keep only E/G/I, and set the new K status.
$ws.Range("A10:D10").ClearContents()
$ws.Range("F10").ClearContents()
$ws.Range("H10").ClearContents()
$ws.Range("J10").ClearContents()
$ws.Range("K10").Value = "FALSO - Documento sólo en Instructores"

# --- Row 11 (new): a record that exists only on the "Sofía" side
# (doc 250188492 / NORA MARIA BLANDON SERNA).
$ws.Range("A11").Value = 3031278
$ws.Range("B11").Value = "COMPLEMENTARIA"
$ws.Range("C11").Value = "CURSO ESPECIAL"
$ws.Range("D11").Value = "ELABORACION DE PRODUCTOS DE REPOSTERIA."
$ws.Range("F11").Value = "CC"
$ws.Range("H11").Value = 250188492
$ws.Range("J11").Value = "NORA MARIA BLANDON SERNA"
$ws.Range("K11").Value = "FALSO - Documento sólo en Sofía"
